# Auto-generated edit script applying cached 'profit' value updates
# produced by the external leve-profit scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 2992.9644
$ws.Range("I76").Value = 2992.9644
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 2992.9644
$ws.Range("L76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -2677.9644
$ws.Range("H79").Value = 2992.9644
$ws.Range("I79").Value = 2992.9644
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 2992.9644
$ws.Range("L79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -1900.9644
$ws.Range("H132").Value = 6192.2666
$ws.Range("I132").Value = 2063.7273
$ws.Range("J132").Value = 17545.75
$ws.Range("K132").Value = 6191.1819
$ws.Range("L132").Value = 52637.25
$ws.Range("M132").Value = -3661.1819
$ws.Range("N132").Value = -57697.25
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H76").Value = 26000
$ws.Range("J76").Value = 26000
$ws.Range("L76").Value = 26000
$ws.Range("N76").Value = -26676
$ws.Range("H79").Value = 26000
$ws.Range("J79").Value = 26000
$ws.Range("L79").Value = 26000
$ws.Range("N79").Value = -28340
$ws.Range("H102").Value = 2921.375
$ws.Range("I102").Value = 2711.6667
$ws.Range("J102").Value = 3550.5
$ws.Range("K102").Value = 2711.6667
$ws.Range("L102").Value = 3550.5
$ws.Range("M102").Value = -1089.6667
$ws.Range("N102").Value = -6794.5
$ws.Range("H122").Value = 1647.091
$ws.Range("I122").Value = 1702
$ws.Range("J122").Value = 1400
$ws.Range("K122").Value = 5106
$ws.Range("L122").Value = 4200
$ws.Range("M122").Value = -2656
$ws.Range("N122").Value = -9100
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H51").Value = 29900
$ws.Range("J51").Value = 29900
$ws.Range("L51").Value = 29900
$ws.Range("N51").Value = -30882
$ws.Range("H55").Value = 29259.334
$ws.Range("J55").Value = 29259.334
$ws.Range("L55").Value = 29259.334
$ws.Range("N55").Value = -29805.334
$ws.Range("H86").Value = 2313.76
$ws.Range("I86").Value = 1779.2
$ws.Range("J86").Value = 3115.6
$ws.Range("K86").Value = 1779.2
$ws.Range("L86").Value = 3115.6
$ws.Range("M86").Value = -656.2
$ws.Range("N86").Value = -5361.6
$ws.Range("H89").Value = 2313.76
$ws.Range("I89").Value = 1779.2
$ws.Range("J89").Value = 3115.6
$ws.Range("K89").Value = 8896
$ws.Range("L89").Value = 15578
$ws.Range("M89").Value = -3280
$ws.Range("N89").Value = -26810
$ws.Range("H99").Value = 2268.6191
$ws.Range("I99").Value = 2070
$ws.Range("J99").Value = 2330.6875
$ws.Range("K99").Value = 2070
$ws.Range("L99").Value = 2330.6875
$ws.Range("M99").Value = -572
$ws.Range("N99").Value = -5326.6875
$ws.Range("H105").Value = 2178.776
$ws.Range("I105").Value = 2207.7144
$ws.Range("J105").Value = 2021.2222
$ws.Range("K105").Value = 2207.7144
$ws.Range("L105").Value = 2021.2222
$ws.Range("M105").Value = -460.7143999999998
$ws.Range("N105").Value = -5515.2222
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 30000
$ws.Range("J68").Value = 30000
$ws.Range("L68").Value = 30000
$ws.Range("N68").Value = -31498
$ws.Range("H71").Value = 30000
$ws.Range("J71").Value = 30000
$ws.Range("L71").Value = 90000
$ws.Range("N71").Value = -97488
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 6074.2856
$ws.Range("I134").Value = 3130
$ws.Range("J134").Value = 10000
$ws.Range("K134").Value = 9390
$ws.Range("L134").Value = 30000
$ws.Range("M134").Value = -4320
$ws.Range("N134").Value = -40140
$ws.Range("H136").Value = 3225.3845
$ws.Range("I136").Value = 1193
$ws.Range("K136").Value = 3579
$ws.Range("M136").Value = 1521
$ws.Range("H138").Value = 1211.3334
$ws.Range("I138").Value = 959.2308
$ws.Range("J138").Value = 2850
$ws.Range("K138").Value = 2877.6924
$ws.Range("L138").Value = 8550
$ws.Range("M138").Value = 2262.3076
$ws.Range("N138").Value = -18830
$ws.Range("H139").Value = 2167.9678
$ws.Range("I139").Value = 850.5833
$ws.Range("J139").Value = 3000
$ws.Range("K139").Value = 2551.7499
$ws.Range("L139").Value = 9000
$ws.Range("M139").Value = 2588.2501
$ws.Range("N139").Value = -19280
$ws.Range("H141").Value = 4008.1191
$ws.Range("I141").Value = 1181.75
$ws.Range("K141").Value = 3545.25
$ws.Range("M141").Value = 1634.75
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H74").Value = 25927.5
$ws.Range("J74").Value = 25927.5
$ws.Range("L74").Value = 25927.5
$ws.Range("N74").Value = -27799.5
$ws.Range("H77").Value = 25927.5
$ws.Range("J77").Value = 25927.5
$ws.Range("L77").Value = 77782.5
$ws.Range("N77").Value = -87142.5
$ws.Range("H80").Value = 2318.3157
$ws.Range("I80").Value = 2304
$ws.Range("J80").Value = 2342.8572
$ws.Range("K80").Value = 2304
$ws.Range("L80").Value = 2342.8572
$ws.Range("M80").Value = -1306
$ws.Range("N80").Value = -4338.8572
$ws.Range("H83").Value = 2318.3157
$ws.Range("I83").Value = 2304
$ws.Range("J83").Value = 2342.8572
$ws.Range("K83").Value = 11520
$ws.Range("L83").Value = 11714.286
$ws.Range("M83").Value = -6528
$ws.Range("N83").Value = -21698.286
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1807.6666
$ws.Range("I40").Value = 1739.375
$ws.Range("J40").Value = 1885.7142
$ws.Range("K40").Value = 1739.375
$ws.Range("L40").Value = 1885.7142
$ws.Range("M40").Value = -1603.375
$ws.Range("N40").Value = -2157.7142
$ws.Range("H68").Value = 10450
$ws.Range("I68").Value = 15742.857
$ws.Range("J68").Value = 3040
$ws.Range("K68").Value = 15742.857
$ws.Range("L68").Value = 3040
$ws.Range("M68").Value = -14993.857
$ws.Range("N68").Value = -4538
$ws.Range("H71").Value = 10450
$ws.Range("I71").Value = 15742.857
$ws.Range("J71").Value = 3040
$ws.Range("K71").Value = 78714.285
$ws.Range("L71").Value = 15200
$ws.Range("M71").Value = -74970.285
$ws.Range("N71").Value = -22688
$ws.Range("H133").Value = 32326
$ws.Range("J133").Value = 32326
$ws.Range("L133").Value = 32326
$ws.Range("N133").Value = -37386
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1032.625
$ws.Range("I81").Value = 894.4286
$ws.Range("J81").Value = 2000
$ws.Range("K81").Value = 1788.8572
$ws.Range("L81").Value = 4000
$ws.Range("M81").Value = -727.8571999999999
$ws.Range("N81").Value = -6122
$ws.Range("H84").Value = 1032.625
$ws.Range("I84").Value = 894.4286
$ws.Range("J84").Value = 2000
$ws.Range("K84").Value = 8944.286
$ws.Range("L84").Value = 20000
$ws.Range("M84").Value = -3640.286
$ws.Range("N84").Value = -30608
$ws.Range("H96").Value = 1565.5333
$ws.Range("J96").Value = 1690.3
$ws.Range("L96").Value = 1690.3
$ws.Range("N96").Value = -4436.3
